$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.55"
$ws.Range("E2").Value = "'3.47%"
$ws.Range("D3").Value = "'40.18"
$ws.Range("E3").Value = "'6.10%"
$ws.Range("D4").Value = "'5.820"
$ws.Range("E4").Value = "'12.23%"
$ws.Range("D5").Value = "'0.08032"
$ws.Range("E5").Value = "'0.55%"
$ws.Range("D6").Value = "'4.576"
$ws.Range("E6").Value = "'1.89%"
$ws.Range("D7").Value = "'8.725"
$ws.Range("E7").Value = "'2.14%"
$ws.Range("D8").Value = "'1.936"
$ws.Range("E8").Value = "'0.44%"
$ws.Range("D10").Value = "'0.9447"
$ws.Range("E10").Value = "'0.14%"
$ws.Range("D11").Value = "'0.1249"
$ws.Range("E11").Value = "'-4.68%"
$ws.Range("D12").Value = "'0.1961"
$ws.Range("E12").Value = "'1.23%"
$ws.Range("D13").Value = "'8.885"
$ws.Range("E13").Value = "'34.59%"
$ws.Range("D14").Value = "'0.09198"
$ws.Range("E14").Value = "'1.49%"
$ws.Range("D15").Value = "'0.03588"
$ws.Range("E15").Value = "'5.66%"
$ws.Range("D16").Value = "'0.09628"
$ws.Range("E16").Value = "'1.12%"
$ws.Range("D17").Value = "'0.001305"
$ws.Range("E17").Value = "'-6.77%"
$ws.Range("D18").Value = "'0.006560"
$ws.Range("E18").Value = "'9.62%"
$ws.Range("D19").Value = "'3.370"
$ws.Range("E19").Value = "'-1.79%"
$ws.Range("D20").Value = "'0.3525"
$ws.Range("E20").Value = "'0.26%"
$ws.Range("D21").Value = "'0.1433"
$ws.Range("E21").Value = "'10.25%"
$ws.Range("D22").Value = "'0.2415"
$ws.Range("E22").Value = "'-0.21%"
$ws.Range("D23").Value = "'0.04397"
$ws.Range("E23").Value = "'0.45%"
$ws.Range("D24").Value = "'0.001261"
$ws.Range("E24").Value = "'2.61%"
$ws.Range("D25").Value = "'0.004320"
$ws.Range("E25").Value = "'1.23%"
$ws.Range("D26").Value = "'0.0001145"
$ws.Range("E26").Value = "'-13.71%"
$ws.Range("E27").Value = "'0.34%"
$ws.Range("D39").Value = "'0.02417"
$ws.Range("E39").Value = "'0.98%"
$ws.Range("D40").Value = "'0.05289"
$ws.Range("E40").Value = "'2.63%"
$ws.Range("D41").Value = "'0.007485"
$ws.Range("E41").Value = "'-1.80%"
$ws.Range("E42").Value = "'1.72%"
$ws.Range("D43").Value = "'0.008764"
$ws.Range("E43").Value = "'2.81%"
$ws.Range("D44").Value = "'0.002107"
$ws.Range("E44").Value = "'0.13%"
$ws.Range("D45").Value = "'0.01056"
$ws.Range("E45").Value = "'21.48%"
$ws.Range("D46").Value = "'0.00006889"
$ws.Range("E46").Value = "'6.59%"
$ws.Range("E47").Value = "'0.63%"
$ws.Range("D48").Value = "'0.003157"
$ws.Range("E48").Value = "'10.40%"
$ws.Range("D49").Value = "'0.001425"
$ws.Range("E49").Value = "'-15.46%"
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("E50").Value = "'0.63%"
$ws.Range("D51").Value = "'0.0002008"
$ws.Range("E51").Value = "'0.63%"
